$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column D (4th col) gets a bit wider (stored width 32 -> 39). The
# ColumnWidth property adds Excel's standard ~0.8333 char padding on write,
# so back it out to land on an exact stored width of 39. ---
$ws.Columns.Item(4).ColumnWidth = 38.166666666666664

# --- Row 2: brand-new listing inserted at the top ---
$ws.Cells.Item(2,1).Value = "2025-12-29 18:28:01"
$ws.Cells.Item(2,2).Value = "建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集"
$ws.Cells.Item(2,3).Value = "システム開発"
$ws.Cells.Item(2,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(2,5).Value = "期限情報なし"
$ws.Cells.Item(2,6).Value = "https://www.lancers.jp/work/detail/5434128"
$ws.Cells.Item(2,7).Value = 368
$ws.Cells.Item(2,8).Value = "🔥AI,Ai ◆開発"

# --- Row 3: previously row 2, date refreshed ---
$ws.Cells.Item(3,1).Value = "2025-12-29 18:28:01"
$ws.Cells.Item(3,2).Value = "大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)"
$ws.Cells.Item(3,3).Value = "システム開発"
$ws.Cells.Item(3,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(3,5).Value = "期限情報なし"
$ws.Cells.Item(3,6).Value = "https://www.lancers.jp/work/detail/5427956"
$ws.Cells.Item(3,7).Value = 310
$ws.Cells.Item(3,8).Value = "🔥AI,Ai"

# --- Row 4: brand-new listing ---
$ws.Cells.Item(4,1).Value = "2025-12-29 18:28:01"
$ws.Cells.Item(4,2).Value = "【急募】製造業向け「製造副産物」の状態(硬度)判定AIのフィジビリティ検証(画像認識/動画解析)"
$ws.Cells.Item(4,3).Value = "システム開発"
$ws.Cells.Item(4,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(4,5).Value = "期限情報なし"
$ws.Cells.Item(4,6).Value = "https://www.lancers.jp/work/detail/5439158"
$ws.Cells.Item(4,7).Value = 303
$ws.Cells.Item(4,8).Value = "🔥AI,Ai"

# --- Row 5: brand-new listing ---
$ws.Cells.Item(5,1).Value = "2025-12-29 18:28:01"
$ws.Cells.Item(5,2).Value = "【急募】宿泊業向けSaaSの予約者取得システム開発"
$ws.Cells.Item(5,3).Value = "システム開発"
$ws.Cells.Item(5,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(5,5).Value = "期限情報なし"
$ws.Cells.Item(5,6).Value = "https://www.lancers.jp/work/detail/5460405"
$ws.Cells.Item(5,7).Value = 118
$ws.Cells.Item(5,8).Value = "◆開発,システム開発"

# --- Row 6: previously row 3, date refreshed ---
$ws.Cells.Item(6,1).Value = "2025-12-29 18:28:01"
$ws.Cells.Item(6,2).Value = "【Unity/XRエンジニア募集】製造業DX支援!既存システムと連携するXRアプリ開発"
$ws.Cells.Item(6,3).Value = "システム開発"
$ws.Cells.Item(6,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(6,5).Value = "期限情報なし"
$ws.Cells.Item(6,6).Value = "https://www.lancers.jp/work/detail/5454210"
$ws.Cells.Item(6,7).Value = 108
$ws.Cells.Item(6,8).Value = "◆開発 ◇アプリ"

# --- Row 7: previously row 4, date refreshed ---
$ws.Cells.Item(7,1).Value = "2025-12-29 18:28:01"
$ws.Cells.Item(7,2).Value = "【SNSアプリ開発】AndroidとiOSのインスタグラム風アプリ制作依頼"
$ws.Cells.Item(7,3).Value = "システム開発"
$ws.Cells.Item(7,4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(7,5).Value = "期限情報なし"
$ws.Cells.Item(7,6).Value = "https://www.lancers.jp/work/detail/5462964"
$ws.Cells.Item(7,7).Value = 100
$ws.Cells.Item(7,8).Value = "◆開発 ◇アプリ"

# --- Row 8: previously row 5, date refreshed and price/count updated ---
$ws.Cells.Item(8,1).Value = "2025-12-29 18:28:01"
$ws.Cells.Item(8,2).Value = "【依頼内容|Googleスプレッドシート構築】"
$ws.Cells.Item(8,3).Value = "システム開発"
$ws.Cells.Item(8,4).Value = "1,000 円 ~ 2,000 円 / 募集期間 7 日、取引期間 0 日"
$ws.Cells.Item(8,5).Value = "期限情報なし"
$ws.Cells.Item(8,6).Value = "https://www.lancers.jp/work/detail/5462891"
$ws.Cells.Item(8,7).Value = 10

# --- Row 9: previously row 6, date refreshed ---
$ws.Cells.Item(9,1).Value = "2025-12-29 18:28:01"
$ws.Cells.Item(9,2).Value = "【急募】VPS(Ubuntu)フロント・バック統合/Docker構築・移行"
$ws.Cells.Item(9,3).Value = "システム開発"
$ws.Cells.Item(9,4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(9,5).Value = "期限情報なし"
$ws.Cells.Item(9,6).Value = "https://www.lancers.jp/work/detail/5462712"
$ws.Cells.Item(9,7).Value = 10

# --- Rebuild the URL hyperlinks from scratch: the old anchors don't track
# the row shuffle above, so clear everything first and re-add one per row
# in the final F2:F9 order (this also keeps hyperlink rIds in row order,
# matching the diff). ---
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(2,6), "https://www.lancers.jp/work/detail/5434128")
$ws.Hyperlinks.Add($ws.Cells.Item(3,6), "https://www.lancers.jp/work/detail/5427956")
$ws.Hyperlinks.Add($ws.Cells.Item(4,6), "https://www.lancers.jp/work/detail/5439158")
$ws.Hyperlinks.Add($ws.Cells.Item(5,6), "https://www.lancers.jp/work/detail/5460405")
$ws.Hyperlinks.Add($ws.Cells.Item(6,6), "https://www.lancers.jp/work/detail/5454210")
$ws.Hyperlinks.Add($ws.Cells.Item(7,6), "https://www.lancers.jp/work/detail/5462964")
$ws.Hyperlinks.Add($ws.Cells.Item(8,6), "https://www.lancers.jp/work/detail/5462891")
$ws.Hyperlinks.Add($ws.Cells.Item(9,6), "https://www.lancers.jp/work/detail/5462712")
